# Generate Report for Handoff
# Replace the old handoff GUID / timestamps with the new ones produced by this handoff run
# (old guid 34f1a8b7-0238-40ae-baf4-bf4be91d966a -> new guid 2d0bc8b8-8423-4709-9e31-1b4533162f81,
#  old xliff hash 3a75ccf4d72451d169912ef69925c305caf529be -> new hash 46d52fa5063a1b07d36b71eae45edf60352e49d1).

$wb = $excel.ActiveWorkbook

$newGuid = "2d0bc8b8-8423-4709-9e31-1b4533162f81"
$newHash = "46d52fa5063a1b07d36b71eae45edf60352e49d1"

$commitSha = "424fdb41ac0034837c28c64da2f91b86a8f6588b"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-28 14:57:36"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md",
    "",
    "",
    "e2e\$newGuid.md"
) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 14:57:32"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 14:57:36"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
) | Out-Null
